$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns E, F, G with new class labels
$ws.Range("E1").Value = "l23.0"
$ws.Range("F1").Value = "l25.0"
$ws.Range("G1").Value = "l27"

# Data rows 2-12: columns E, F, G duplicate the values of columns B, C, D
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 4).Value()
}

# Update the active selection as recorded in the saved file
$ws.Range("E2").Select()
